$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Venkatesh Iyer"

# Force numeric-looking columns (runs, balls, fours, sixes, sr) to be stored as text,
# matching the source data which represents them as strings.
$ws.Range("E2:I11").NumberFormat = "@"
# "states" is blank for match 31st (row 7); force text format so the empty cell persists.
$ws.Range("D7").NumberFormat = "@"

# Header row
$ws.Range("A1").Value = "matchNo"
$ws.Range("B1").Value = "teamName"
$ws.Range("C1").Value = "batterName"
$ws.Range("D1").Value = "states"
$ws.Range("E1").Value = "runs"
$ws.Range("F1").Value = "balls"
$ws.Range("G1").Value = "fours"
$ws.Range("H1").Value = "sixes"
$ws.Range("I1").Value = "sr"
$ws.Range("J1").Value = "opponentTeamName"
$ws.Range("K1").Value = "venue"
$ws.Range("L1").Value = "date"
$ws.Range("M1").Value = "result"

# Data rows
# Row 2
$ws.Range("A2").Value = "Qualifier"
$ws.Range("B2").Value = "Kolkata Knight Riders"
$ws.Range("C2").Value = "Venkatesh Iyer"
$ws.Range("D2").Value = "c sub (SPD Smith) b Rabada"
$ws.Range("E2").Value = "55"
$ws.Range("F2").Value = "41"
$ws.Range("G2").Value = "4"
$ws.Range("H2").Value = "3"
$ws.Range("I2").Value = "134.14"
$ws.Range("J2").Value = "Delhi Capitals"
$ws.Range("K2").Value = "Sharjah"
$ws.Range("L2").Value = "October 13"
$ws.Range("M2").Value = "KKR won by 3 wickets (with 1 ball remaining)"

# Row 3
$ws.Range("A3").Value = "Eliminator"
$ws.Range("B3").Value = "Kolkata Knight Riders"
$ws.Range("C3").Value = "Venkatesh Iyer"
$ws.Range("D3").Value = "c †Bharat b Patel"
$ws.Range("E3").Value = "26"
$ws.Range("F3").Value = "30"
$ws.Range("G3").Value = "0"
$ws.Range("H3").Value = "1"
$ws.Range("I3").Value = "86.66"
$ws.Range("J3").Value = "Royal Challengers Bangalore"
$ws.Range("K3").Value = "Sharjah"
$ws.Range("L3").Value = "October 11"
$ws.Range("M3").Value = "KKR won by 4 wickets (with 2 balls remaining)"

# Row 4
$ws.Range("A4").Value = "41st"
$ws.Range("B4").Value = "Kolkata Knight Riders"
$ws.Range("C4").Value = "Venkatesh Iyer"
$ws.Range("D4").Value = "b Lalit Yadav"
$ws.Range("E4").Value = "14"
$ws.Range("F4").Value = "15"
$ws.Range("G4").Value = "2"
$ws.Range("H4").Value = "0"
$ws.Range("I4").Value = "93.33"
$ws.Range("J4").Value = "Delhi Capitals"
$ws.Range("K4").Value = "Sharjah"
$ws.Range("L4").Value = "September 28"
$ws.Range("M4").Value = "KKR won by 3 wickets (with 10 balls remaining)"

# Row 5
$ws.Range("A5").Value = "Final"
$ws.Range("B5").Value = "Kolkata Knight Riders"
$ws.Range("C5").Value = "Venkatesh Iyer"
$ws.Range("D5").Value = "c Jadeja b Thakur"
$ws.Range("E5").Value = "50"
$ws.Range("F5").Value = "32"
$ws.Range("G5").Value = "5"
$ws.Range("H5").Value = "3"
$ws.Range("I5").Value = "156.25"
$ws.Range("J5").Value = "Chennai Super Kings"
$ws.Range("K5").Value = "Dubai (DSC)"
$ws.Range("L5").Value = "October 15"
$ws.Range("M5").Value = "Super Kings won by 27 runs"

# Row 6
$ws.Range("A6").Value = "45th"
$ws.Range("B6").Value = "Kolkata Knight Riders"
$ws.Range("C6").Value = "Venkatesh Iyer"
$ws.Range("D6").Value = "c Hooda b Ravi Bishnoi"
$ws.Range("E6").Value = "67"
$ws.Range("F6").Value = "49"
$ws.Range("G6").Value = "9"
$ws.Range("H6").Value = "1"
$ws.Range("I6").Value = "136.73"
$ws.Range("J6").Value = "Punjab Kings"
$ws.Range("K6").Value = "Dubai (DSC)"
$ws.Range("L6").Value = "October 01"
$ws.Range("M6").Value = "Punjab Kings won by 5 wickets (with 3 balls remaining)"

# Row 7
$ws.Range("A7").Value = "31st"
$ws.Range("B7").Value = "Kolkata Knight Riders"
$ws.Range("C7").Value = "Venkatesh Iyer"
$ws.Range("D7").Formula = ""
$ws.Range("E7").Value = "41"
$ws.Range("F7").Value = "27"
$ws.Range("G7").Value = "7"
$ws.Range("H7").Value = "1"
$ws.Range("I7").Value = "151.85"
$ws.Range("J7").Value = "Royal Challengers Bangalore"
$ws.Range("K7").Value = "Abu Dhabi"
$ws.Range("L7").Value = "September 20"
$ws.Range("M7").Value = "KKR won by 9 wickets (with 60 balls remaining)"

# Row 8
$ws.Range("A8").Value = "49th"
$ws.Range("B8").Value = "Kolkata Knight Riders"
$ws.Range("C8").Value = "Venkatesh Iyer"
$ws.Range("D8").Value = "c Williamson b Holder"
$ws.Range("E8").Value = "8"
$ws.Range("F8").Value = "14"
$ws.Range("G8").Value = "0"
$ws.Range("H8").Value = "0"
$ws.Range("I8").Value = "57.14"
$ws.Range("J8").Value = "Sunrisers Hyderabad"
$ws.Range("K8").Value = "Dubai (DSC)"
$ws.Range("L8").Value = "October 03"
$ws.Range("M8").Value = "KKR won by 6 wickets (with 2 balls remaining)"

# Row 9
$ws.Range("A9").Value = "54th"
$ws.Range("B9").Value = "Kolkata Knight Riders"
$ws.Range("C9").Value = "Venkatesh Iyer"
$ws.Range("D9").Value = "b Tewatia"
$ws.Range("E9").Value = "38"
$ws.Range("F9").Value = "35"
$ws.Range("G9").Value = "3"
$ws.Range("H9").Value = "2"
$ws.Range("I9").Value = "108.57"
$ws.Range("J9").Value = "Rajasthan Royals"
$ws.Range("K9").Value = "Sharjah"
$ws.Range("L9").Value = "October 07"
$ws.Range("M9").Value = "KKR won by 86 runs"

# Row 10
$ws.Range("A10").Value = "38th"
$ws.Range("B10").Value = "Kolkata Knight Riders"
$ws.Range("C10").Value = "Venkatesh Iyer"
$ws.Range("D10").Value = "c †Dhoni b Thakur"
$ws.Range("E10").Value = "18"
$ws.Range("F10").Value = "15"
$ws.Range("G10").Value = "3"
$ws.Range("H10").Value = "0"
$ws.Range("I10").Value = "120.00"
$ws.Range("J10").Value = "Chennai Super Kings"
$ws.Range("K10").Value = "Abu Dhabi"
$ws.Range("L10").Value = "September 26"
$ws.Range("M10").Value = "Super Kings won by 2 wickets"

# Row 11
$ws.Range("A11").Value = "34th"
$ws.Range("B11").Value = "Kolkata Knight Riders"
$ws.Range("C11").Value = "Venkatesh Iyer"
$ws.Range("D11").Value = "b Bumrah"
$ws.Range("E11").Value = "53"
$ws.Range("F11").Value = "30"
$ws.Range("G11").Value = "4"
$ws.Range("H11").Value = "3"
$ws.Range("I11").Value = "176.66"
$ws.Range("J11").Value = "Mumbai Indians"
$ws.Range("K11").Value = "Abu Dhabi"
$ws.Range("L11").Value = "September 23"
$ws.Range("M11").Value = "KKR won by 7 wickets (with 29 balls remaining)"

